$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 285, shifting existing rows 285-307 down to 287-309
$ws.Rows.Item(285).Insert()
$ws.Rows.Item(285).Insert()

# Row 285
$ws.Cells.Item(285, 1).Value = 7
$ws.Cells.Item(285, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(285, 3).Value = 'Ñuble'
$ws.Cells.Item(285, 4).Value = 44931
$ws.Cells.Item(285, 5).Value = 16
$ws.Cells.Item(285, 6).Value = 'Fruta'
$ws.Cells.Item(285, 7).Value = 100103
$ws.Cells.Item(285, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(285, 9).Value = 100103004
$ws.Cells.Item(285, 10).Value = 'Durazno'
$ws.Cells.Item(285, 11).Value = 'Kurakata'
$ws.Cells.Item(285, 12).Value = 'Primera'
$ws.Cells.Item(285, 13).Value = 120
$ws.Cells.Item(285, 14).Value = 15000
$ws.Cells.Item(285, 15).Value = 16000
$ws.Cells.Item(285, 16).Value = 15500
$ws.Cells.Item(285, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(285, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(285, 19).Value = 1033
$ws.Cells.Item(285, 20).Value = 15

# Row 286
$ws.Cells.Item(286, 1).Value = 7
$ws.Cells.Item(286, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(286, 3).Value = 'Ñuble'
$ws.Cells.Item(286, 4).Value = 44931
$ws.Cells.Item(286, 5).Value = 16
$ws.Cells.Item(286, 6).Value = 'Fruta'
$ws.Cells.Item(286, 7).Value = 100103
$ws.Cells.Item(286, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(286, 9).Value = 100103004
$ws.Cells.Item(286, 10).Value = 'Durazno'
$ws.Cells.Item(286, 11).Value = 'Kurakata'
$ws.Cells.Item(286, 12).Value = 'Segunda'
$ws.Cells.Item(286, 13).Value = 100
$ws.Cells.Item(286, 14).Value = 14000
$ws.Cells.Item(286, 15).Value = 14000
$ws.Cells.Item(286, 16).Value = 14000
$ws.Cells.Item(286, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(286, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(286, 19).Value = 933
$ws.Cells.Item(286, 20).Value = 15
